$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = -7.328999999999999
$ws.Range("D32").Value = -8.01
$ws.Range("D36").Value = -8.051
$ws.Range("D38").Value = -7.662000000000001
$ws.Range("D46").Value = -8.190000000000001
$ws.Range("D54").Value = -8.494
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("D67").Value = -7.281000000000001
$ws.Range("D69").Value = -7.321000000000001
$ws.Range("D72").Value = -7.434
$ws.Range("D91").Value = -6.970000000000002
$ws.Range("D99").Value = -8.019
